# Apply "updated with new names as of Dec 21 am" edit.
# Replace the full data body (rows 2..N) with the new, expanded roster.
# Columns: A=agency, B=name, C=team_lead (unused/blank), D=on_multiple_teams,
#          E=most_recent_employment, F=source_of_funding

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data set (in final row order, rows 2 through 22)
$data = @(
  @("Department of Education", "David Holmes", "", "Y", "Rebellion Defense", "Volunteer"),
  @("Department of Energy", "Caroline Grey", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Department of Health and Human Services", "Kevin O’Connor", "", "Y", "GW Medical Faculty Associates", "Volunteer"),
  @("Department of Health and Human Services", "Clara Pratte", "", "Y", "Strongbow Strategies", "Transition — PT Fund, Inc."),
  @("Department of Health and Human Services", "Rochelle Walensky", "", "", "Massachusetts General Hospital", "Volunteer"),
  @("Department of Homeland Security", "John Bivona", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Department of Homeland Security", "Kevin Munoz", "", "", "Latino Decisions", "Volunteer"),
  @("Department of Housing and Urban Development", "Analysse Escobar", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Department of Housing and Urban Development", "Mikayla Ferrell", "", "", "Pennsylvania Democratic Party", "Transition — PT Fund, Inc."),
  @("Department of Justice", "Theresa Bradley", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Department of the Interior", "Maggie Thompson", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Department of the Treasury", "William Doerrer", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Department of Transportation", "Allie Panther", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Environmental Protection Agency", "Sinceré Harris", "", "", "Pennsylvania Democratic Party", "Transition — PT Fund, Inc."),
  @("Executive Office of the President, Management and Administration", "Rob Flaherty", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Executive Office of the President, Management and Administration", "Dan Jacobson", "", "", "Arnold & Porter Kaye Scholer LLP", "Volunteer"),
  @("Executive Office of the President, Management and Administration", "Kevin O’Connor", "", "Y", "GW Medical Faculty Associates", "Volunteer"),
  @("Intelligence Community", "Avril Haines", "", "", "Columbia University", "Transition — PT Fund, Inc."),
  @("International Development", "Zeppa Kreager", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("Office of Personnel Management", "Jason Tengco", "", "", "Biden for President", "Transition — PT Fund, Inc."),
  @("United States Department of Agriculture", "Lexi Coburn", "", "", "North Carolina Democratic Party", "Transition — PT Fund, Inc.")
)

$newLastRow = 1 + $data.Count

# Clear out any existing data rows below the header before rewriting, in case
# the old sheet had more rows than the new one (not the case here, but safe).
$oldUsedRows = $ws.UsedRange.Rows.Count
if ($oldUsedRows -gt $newLastRow) {
  $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldUsedRows, 6))
  $clearRange.ClearContents()
}

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = 2 + $i
  $row = $data[$i]
  for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}

$ws.Range("A1:F$newLastRow").EntireColumn.AutoFit() | Out-Null
